# Update getSeconds to add resolution even if HR clock rolls over
#
# The only authored input in this sheet is C16 (named range "ExpectedOutputS").
# Every other changed cell (C18, C19, C23, C24, C29-C32, C36-C39, C43-C46) is a
# formula that depends on ExpectedOutputS, either directly or transitively, so
# updating C16 and letting Excel recalculate reproduces all of the downstream
# value changes in one step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ExpectedOutputS (named range -> Sheet1!$C$16)
$ws.Range("ExpectedOutputS").Value = 5.3453429999999997

# Match the author's final selection/scroll position recorded in the sheet view.
$ws.Range("C14").Select() | Out-Null
